$d = $word.ActiveDocument
# This script applies 4 targeted edits to the diary (Dagbok) document,
# each via InsertXML'ing a reconstructed paragraph so run/proofErr
# boundaries come out exactly as in the target revision.

# --- Edit 1: paragraph 2 ---
$xml1 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4FDFC775" w14:textId="47F21343" w:rsidR="00734070" w:rsidRDefault="0032762E" w:rsidP="0032762E"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r w:rsidRPr="006407C0"><w:rPr><w:rStyle w:val="Rubrik1Char"/></w:rPr><w:t>2023-01-09</w:t></w:r><w:r w:rsidRPr="006407C0"><w:rPr><w:rStyle w:val="Rubrik1Char"/></w:rPr><w:br/></w:r><w:r><w:t xml:space="preserve">Tänkt ut hur </w:t></w:r><w:r w:rsidR="00B362BF"><w:t>undermenyerna ska se ut</w:t></w:r><w:r><w:t xml:space="preserve"> och skr</w:t></w:r><w:r w:rsidR="00511886"><w:t>ivit</w:t></w:r><w:r><w:t xml:space="preserve"> in detta som en kommentar. </w:t></w:r><w:r w:rsidR="00481411"><w:br/><w:t>Gjor</w:t></w:r><w:r w:rsidR="00511886"><w:t>t</w:t></w:r><w:r w:rsidR="00481411"><w:t xml:space="preserve"> en koppling till databas via </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00481411"><w:t>Code</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00481411"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00481411"><w:t>First</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00354D20"><w:t xml:space="preserve"> och skapa</w:t></w:r><w:r w:rsidR="00511886"><w:t>t</w:t></w:r><w:r w:rsidR="00354D20"><w:t xml:space="preserve"> en databas från koden. </w:t></w:r><w:r w:rsidR="00734070"><w:br/><w:t>Skapa</w:t></w:r><w:r w:rsidR="001935BA"><w:t>t</w:t></w:r><w:r w:rsidR="00630774"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00734070"><w:t xml:space="preserve">ett </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00734070"><w:t>class</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00734070"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00734070"><w:t>library</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00734070"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00630774"><w:t>för beräkningar. Gjor</w:t></w:r><w:r w:rsidR="001935BA"><w:t>t</w:t></w:r><w:r w:rsidR="00630774"><w:t xml:space="preserve"> ett till </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00630774"><w:t>class</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00630774"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00630774"><w:t>library</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00630774"><w:t xml:space="preserve"> för felhantering. Funderade på om det skulle vara ett och samma </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00630774"><w:t>class</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00630774"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00630774"><w:t>library</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00630774"><w:t xml:space="preserve"> men </w:t></w:r><w:r w:rsidR="001423DB"><w:t>t</w:t></w:r><w:r w:rsidR="00E63664"><w:t>än</w:t></w:r><w:r w:rsidR="001423DB"><w:t xml:space="preserve">kte </w:t></w:r><w:r w:rsidR="00630774"><w:t xml:space="preserve">att det blir enklare att hantera dem ifall </w:t></w:r><w:r w:rsidR="009452F9"><w:t xml:space="preserve">varje </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="009452F9"><w:t>class</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="009452F9"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="009452F9"><w:t>library</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="009452F9"><w:t xml:space="preserve"> bara innehåller en sak.</w:t></w:r><w:r w:rsidR="000574BC"><w:t xml:space="preserve"> Om jag skulle vilja använda exempelvis </w:t></w:r><w:r w:rsidR="001423DB"><w:t xml:space="preserve">felhantering </w:t></w:r><w:r w:rsidR="000574BC"><w:t>i ett annat projekt så behöver jag inte</w:t></w:r><w:r w:rsidR="001423DB"><w:t xml:space="preserve"> ha beräkningar med i samma </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001423DB"><w:t>class</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001423DB"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001423DB"><w:t>library</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001423DB"><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidR="000574BC"><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="256">
<pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/styles" Target="styles.xml"/></Relationships></pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/styles.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.styles+xml">
<pkg:xmlData><w:styles xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:style w:type="paragraph" w:default="1" w:styleId="Normal"><w:name w:val="Normal"/></w:style><w:style w:type="character" w:customStyle="1" w:styleId="Rubrik1Char"><w:name w:val="Rubrik 1 Char"/></w:style></w:styles></pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$word.ActiveDocument.Paragraphs.Item(2).Range.InsertXML($xml1)

# --- Edit 2: paragraph 5 ---
$xml2 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2329CB71" w14:textId="3F701BF6" w:rsidR="00BB169E" w:rsidRDefault="00BB169E" w:rsidP="0032762E"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r w:rsidRPr="00BB169E"><w:rPr><w:rStyle w:val="Rubrik1Char"/></w:rPr><w:t>2023-01-12</w:t></w:r><w:r w:rsidRPr="00BB169E"><w:rPr><w:rStyle w:val="Rubrik1Char"/></w:rPr><w:br/></w:r><w:r><w:t xml:space="preserve">Försökt få till </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>strategy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pattern</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> i </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Calculate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> för att slippa upprepa kod. Hade visst redan börjat med </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>strategy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pattern</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> i </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>shapes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> så jag fortsatte med det också.  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="256">
<pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/styles" Target="styles.xml"/></Relationships></pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/styles.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.styles+xml">
<pkg:xmlData><w:styles xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:style w:type="paragraph" w:default="1" w:styleId="Normal"><w:name w:val="Normal"/></w:style><w:style w:type="character" w:customStyle="1" w:styleId="Rubrik1Char"><w:name w:val="Rubrik 1 Char"/></w:style></w:styles></pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$word.ActiveDocument.Paragraphs.Item(5).Range.InsertXML($xml2)

# --- Edit 3: paragraph 8 ---
$xml3 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="152EFC66" w14:textId="77777777" w:rsidR="00316B37" w:rsidRDefault="00316B37" w:rsidP="00316B37"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>ATT GÖRA:</w:t></w:r><w:r><w:br/><w:t xml:space="preserve"># Miniräknare - redovisa svar med två decimaler </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>vb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="256">
<pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/styles" Target="styles.xml"/></Relationships></pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/styles.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.styles+xml">
<pkg:xmlData><w:styles xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:style w:type="paragraph" w:default="1" w:styleId="Normal"><w:name w:val="Normal"/></w:style><w:style w:type="character" w:customStyle="1" w:styleId="Rubrik1Char"><w:name w:val="Rubrik 1 Char"/></w:style></w:styles></pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$word.ActiveDocument.Paragraphs.Item(8).Range.InsertXML($xml3)

# --- Edit 4: paragraph 14 ---
$xml4 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7E5D3E13" w14:textId="452CDD01" w:rsidR="00316B37" w:rsidRDefault="00316B37" w:rsidP="00316B37"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>?? Ska användaren mata in två tal även vid roten ur ??</w:t></w:r></w:p></w:body></w:document></pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="256">
<pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/styles" Target="styles.xml"/></Relationships></pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/styles.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.styles+xml">
<pkg:xmlData><w:styles xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:style w:type="paragraph" w:default="1" w:styleId="Normal"><w:name w:val="Normal"/></w:style><w:style w:type="character" w:customStyle="1" w:styleId="Rubrik1Char"><w:name w:val="Rubrik 1 Char"/></w:style></w:styles></pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$word.ActiveDocument.Paragraphs.Item(14).Range.InsertXML($xml4)

